{"js": "// Convert the two plain-text \"\u25a1\" (U+25A1) checkbox glyphs in the\n// \"Nella sua qualit\u00e0 di:\" table (rows \"Titolare o Legale rappresentante\"\n// and \"Procuratore\") into real Word legacy checkbox content controls\n// (w:sdt / w14:checkbox), matching Word's own behaviour when you use\n// Developer \u25b8 Check Box Content Control: the glyph becomes \"\u2610\" (U+2610)\n// rendered in MS Gothic, and the run is wrapped by a checkbox SDT\n// (unchecked, checkedState U+2612, uncheckedState U+2610).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row 3 (0-based) -> \"Titolare o Legale rappresentante\"\n// Row 4 (0-based) -> \"Procuratore\"\nconst checkboxRows = [3, 4];\n\nfor (const rowIndex of checkboxRows) {\n  const cell = table.getCell(rowIndex, 0);\n  const cellRange = cell.body.getRange();\n\n  const checkBoxControl = cellRange.insertContentControl(\"CheckBox\");\n  checkBoxControl.checked = false;\n}\n\nawait context.sync();\n", "ps1": "# Convert the two plain-text \"[]\" (U+25A1, \"\u25a1\") checkbox glyphs in the\n# \"Nella sua qualita di:\" table (rows \"Titolare o Legale rappresentante\"\n# and \"Procuratore\") into real Word legacy checkbox content controls\n# (w:sdt / w14:checkbox), matching Word's own behaviour when you use\n# Developer > Check Box Content Control: the glyph becomes \"[]\" (U+2610,\n# \"\u2610\") rendered in MS Gothic, and the run is wrapped by a checkbox SDT\n# (unchecked, checkedState U+2612, uncheckedState U+2610).\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Column 1 of table rows 4 and 5 hold the checkbox glyph, next to\n# \"Titolare o Legale rappresentante\" and \"Procuratore\" respectively.\n$checkboxRows = 4, 5\n\nforeach ($rowIndex in $checkboxRows) {\n    $cell = $tbl.Cell($rowIndex, 1)\n    $cellRange = $cell.Range\n\n    # Build a fresh Range from explicit character offsets (rather than\n    # handing the Cell's own Range object straight to ContentControls.Add)\n    # and trim off the trailing end-of-cell mark so only the checkbox\n    # character itself gets wrapped.\n    $startPos = $cellRange.Start\n    $endPos = $cellRange.End - 1\n    $checkboxRange = $d.Range($startPos, $endPos)\n\n    $cc = $d.ContentControls.Add(8, $checkboxRange)  # 8 = wdContentControlCheckBox\n    $cc.Checked = $false\n}\n"}
